$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add column A "type" values for rows 2-18
$rowsType1 = 2..11
foreach ($r in $rowsType1) {
    $ws.Cells.Item($r, 1).Value = 1
}

$rowsType3 = 12..18
foreach ($r in $rowsType3) {
    $ws.Cells.Item($r, 1).Value = 3
}

# Set column F (slugDlg) to "no_pic" for rows 2-18
foreach ($r in 2..18) {
    $ws.Cells.Item($r, 6).Value = "no_pic"
}
